$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = [DateTime]::FromOADate(45206)

for ($r = 2; $r -le 507; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
